# Update fuel metrics data (P&L calculation refresh) on Sheet1, rows 2-26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; A = "FL6679"; B = "AC7432"; C = 850; D = 2.03; E = 45839; F = 3.409; G = 5882.23 },
    @{ Row = 3; A = "FL1976"; B = "AC7432"; C = 850; D = 3.68; E = 45809; F = 2.427; G = 7591.66 },
    @{ Row = 4; A = "FL4876"; B = "AC7432"; C = 850; D = 4.95; E = 45809; F = 2.427; G = 10211.6 },
    @{ Row = 5; A = "FL1198"; B = "AC7432"; C = 850; D = 2.03; E = 45809; F = 2.427; G = 4187.79 },
    @{ Row = 6; A = "FL3066"; B = "AC7432"; C = 850; D = 2.5; E = 45809; F = 2.427; G = 5157.38 },
    @{ Row = 7; A = "FL5388"; B = "AC7432"; C = 850; D = 3.68; E = 45809; F = 2.427; G = 7591.66 },
    @{ Row = 8; A = "FL5301"; B = "AC9250"; C = 1000; D = 2.5; E = 45809; F = 2.427; G = 6067.5 },
    @{ Row = 9; A = "FL7875"; B = "AC9250"; C = 1000; D = 4.95; E = 45809; F = 2.427; G = 12013.65 },
    @{ Row = 10; A = "FL6369"; B = "AC9250"; C = 1000; D = 3.68; E = 45809; F = 2.427; G = 8931.360000000001 },
    @{ Row = 11; A = "FL7466"; B = "AC9250"; C = 1000; D = 3.68; E = 45809; F = 2.427; G = 8931.360000000001 },
    @{ Row = 12; A = "FL3108"; B = "AC9250"; C = 1000; D = 2.5; E = 45809; F = 2.427; G = 6067.5 },
    @{ Row = 13; A = "FL4736"; B = "AC9250"; C = 1000; D = 1.45; E = 45809; F = 2.427; G = 3519.15 },
    @{ Row = 14; A = "FL5961"; B = "AC4709"; C = 1350; D = 7.25; E = 45809; F = 2.427; G = 23754.26 },
    @{ Row = 15; A = "FL2734"; B = "AC4709"; C = 1350; D = 7.25; E = 45809; F = 2.427; G = 23754.26 },
    @{ Row = 16; A = "FL1013"; B = "AC4709"; C = 1350; D = 7.25; E = 45809; F = 2.427; G = 23754.26 },
    @{ Row = 17; A = "FL5705"; B = "AC7421"; C = 890; D = 2.03; E = 45839; F = 3.409; G = 6159.04 },
    @{ Row = 18; A = "FL1524"; B = "AC7421"; C = 890; D = 3.68; E = 45809; F = 2.427; G = 7948.91 },
    @{ Row = 19; A = "FL5421"; B = "AC7421"; C = 890; D = 4.95; E = 45809; F = 2.427; G = 10692.15 },
    @{ Row = 20; A = "FL6133"; B = "AC7421"; C = 890; D = 4.95; E = 45809; F = 2.427; G = 10692.15 },
    @{ Row = 21; A = "FL8299"; B = "AC7979"; C = 500; D = 1.45; E = 45809; F = 2.427; G = 1759.58 },
    @{ Row = 22; A = "FL5785"; B = "AC5682"; C = 790; D = 1.45; E = 45839; F = 3.409; G = 3905.01 },
    @{ Row = 23; A = "FL5235"; B = "AC5682"; C = 790; D = 4.95; E = 45809; F = 2.427; G = 9490.780000000001 },
    @{ Row = 24; A = "FL5689"; B = "AC5682"; C = 790; D = 2.03; E = 45809; F = 2.427; G = 3892.18 },
    @{ Row = 25; A = "FL8121"; B = "AC5682"; C = 790; D = 1.45; E = 45809; F = 2.427; G = 2780.13 },
    @{ Row = 26; A = "FL1990"; B = "AC4686"; C = 1350; D = 7.25; E = 45809; F = 2.427; G = 23754.26 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
